$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.881.15'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -0.95%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.854.85'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -1.47%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '598.59'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.92%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '166.51'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.90%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.855.36'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -1.36%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("E9").Value = '  -0.57%  '
$ws.Range("E10").Value = '  -0.30%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.32'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -1.37%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.456'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.57%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000249'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +1.57%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.81'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.11%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.499.67'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -1.41%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.854.41'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -2.73%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '67.852.13'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -1.30%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.17'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +7.05%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.37'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.41%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.111'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -1.90%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.82'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -3.36%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '465.44'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -3.67%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.728'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +1.53%  '
$ws.Range("E24").Value = '  -3.87%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.26'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -1.29%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.22'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.80%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.13'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +1.09%  '
$ws.Range("E29").Value = '  -0.91%  '
$ws.Range("E30").Value = '  +0.12%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.003.67'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -1.43%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.71'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -1.66%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.32'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -2.08%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '31.02'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -2.87%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.829.75'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -0.63%  '
$ws.Range("E36").Value = '  -2.30%  '
$ws.Range("E37").Value = '  +0.16%  '
$ws.Range("E38").Value = '  -2.11%  '
$ws.Range("B39").Value = 'Filecoin'
$ws.Range("C39").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.89'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +0.46%  '
$ws.Range("B40").Value = 'dogwifhat'
$ws.Range("C40").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.26'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +8.65%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.999'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -0.03%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.312'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -1.89%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '428.91'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -1.14%  '
$ws.Range("E44").Value = '  -0.09%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '47.19'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -2.66%  '
$ws.Range("E47").Value = '  +1.17%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '143.55'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +1.22%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '26.31'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +1.19%  '
$ws.Range("B50").Value = 'FLOKI'
$ws.Range("C50").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.000270'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +6.47%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '39.43'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.77%  '
